# Update EIA Table 4.4 workbook: add November 2016 monthly row, update the
# "Year to Date" and "Rolling 12 Months" summary rows, and refresh the
# title / rolling-window label text (October 2016 -> November 2016).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the title and the "Rolling 12 Months Ending in ..." label
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Table 4.4. Receipts, Average Cost, and Quality of Fossil Fuels: Commercial Sector, 2006 - November 2016"

# ------------------------------------------------------------------
# 2. Insert a new row for "November" right above the existing
#    "Year to Date" summary row (old row 53), pushing everything
#    below it down by one row. Use PasteSpecial(formats) from the
#    row directly above (the "October" row) so the new row reuses
#    the existing style indices instead of Excel fabricating new
#    ones.
# ------------------------------------------------------------------
$ws.Rows("53:53").Insert()
$ws.Range("A52:M52").Copy()
$ws.Range("A53:M53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the November 2016 monthly data (row 53)
$ws.Range("A53").Value2 = "November"
$ws.Range("B53").Value2 = 237
$ws.Range("C53").Value2 = 10
$ws.Range("D53").Value2 = "W"
$ws.Range("E53").Value2 = "W"
$ws.Range("F53").Value2 = 3.04
$ws.Range("G53").Value2 = 17.1
$ws.Range("H53").Value2 = 0
$ws.Range("I53").Value2 = 0
$ws.Range("J53").Value2 = "--"
$ws.Range("K53").Value2 = "--"
$ws.Range("L53").Value2 = "--"
$ws.Range("M53").Value2 = 0

# ------------------------------------------------------------------
# 3. Update the "Year to Date" block (now rows 55-57: 2014, 2015, 2016)
# ------------------------------------------------------------------
# 2014
$ws.Range("A55").Value2 = 2014
$ws.Range("B55").Value2 = 3832
$ws.Range("C55").Value2 = 170
$ws.Range("D55").Value2 = "W"
$ws.Range("E55").Value2 = "W"
$ws.Range("F55").Value2 = 2.47
$ws.Range("G55").Value2 = 17.5
$ws.Range("H55").Value2 = 0
$ws.Range("I55").Value2 = 0
$ws.Range("J55").Value2 = "--"
$ws.Range("K55").Value2 = "--"
$ws.Range("L55").Value2 = "--"
$ws.Range("M55").Value2 = 0

# 2015
$ws.Range("A56").Value2 = 2015
$ws.Range("B56").Value2 = 2252
$ws.Range("C56").Value2 = 100
$ws.Range("D56").Value2 = 2.86
$ws.Range("E56").Value2 = 64.22
$ws.Range("F56").Value2 = 2.52
$ws.Range("G56").Value2 = 13.8
$ws.Range("H56").Value2 = 0
$ws.Range("I56").Value2 = 0
$ws.Range("J56").Value2 = "--"
$ws.Range("K56").Value2 = "--"
$ws.Range("L56").Value2 = "--"
$ws.Range("M56").Value2 = 0

# 2016
$ws.Range("A57").Value2 = 2016
$ws.Range("B57").Value2 = 1074
$ws.Range("C57").Value2 = 47
$ws.Range("D57").Value2 = "W"
$ws.Range("E57").Value2 = "W"
$ws.Range("F57").Value2 = 3.03
$ws.Range("G57").Value2 = 7.6
$ws.Range("H57").Value2 = 0
$ws.Range("I57").Value2 = 0
$ws.Range("J57").Value2 = "--"
$ws.Range("K57").Value2 = "--"
$ws.Range("L57").Value2 = "--"
$ws.Range("M57").Value2 = 0

# ------------------------------------------------------------------
# 4. Update the "Rolling 12 Months Ending in ..." header text (now row 58)
# ------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# ------------------------------------------------------------------
# 5. Update the "Rolling 12 Months" block (now rows 59-60: 2015, 2016)
# ------------------------------------------------------------------
# 2015
$ws.Range("A59").Value2 = 2015
$ws.Range("B59").Value2 = 2516
$ws.Range("C59").Value2 = 112
$ws.Range("D59").Value2 = "W"
$ws.Range("E59").Value2 = "W"
$ws.Range("F59").Value2 = 2.55
$ws.Range("G59").Value2 = 13.7
$ws.Range("H59").Value2 = 0
$ws.Range("I59").Value2 = 0
$ws.Range("J59").Value2 = "--"
$ws.Range("K59").Value2 = "--"
$ws.Range("L59").Value2 = "--"
$ws.Range("M59").Value2 = 0

# 2016
$ws.Range("A60").Value2 = 2016
$ws.Range("B60").Value2 = 1262
$ws.Range("C60").Value2 = 56
$ws.Range("D60").Value2 = "W"
$ws.Range("E60").Value2 = "W"
$ws.Range("F60").Value2 = 3
$ws.Range("G60").Value2 = 8
$ws.Range("H60").Value2 = 0
$ws.Range("I60").Value2 = 0
$ws.Range("J60").Value2 = "--"
$ws.Range("K60").Value2 = "--"
$ws.Range("L60").Value2 = "--"
$ws.Range("M60").Value2 = 0
